# HIST.xlsx: split the old single "Terms Typically Offered" column (D) into
# four columns - Corequisites (D), Concurrent (E), Recommended (F) and
# Terms Typically Offered (G) - to reflect the new requirement separation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: insert the three new headers before the existing one and
#    move "Terms Typically Offered" out to column G.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,4).Value = "Corequisites"
$ws.Cells.Item(1,5).Value = "Concurrent"
$ws.Cells.Item(1,6).Value = "Recommended"
$ws.Cells.Item(1,7).Value = "Terms Typically Offered"

# ---------------------------------------------------------------------------
# 2. Data rows (2-102): the old column D value ("Terms Typically Offered")
#    moves to column G, while the new D/E/F default to "NA".
# ---------------------------------------------------------------------------
for ($row = 2; $row -le 102; $row++) {
    $oldTerms = $ws.Cells.Item($row, 4).Value2
    $ws.Cells.Item($row, 4).Value = "NA"
    $ws.Cells.Item($row, 5).Value = "NA"
    $ws.Cells.Item($row, 6).Value = "NA"
    $ws.Cells.Item($row, 7).Value = $oldTerms
}

# ---------------------------------------------------------------------------
# 3. Row-specific corrections, per the source-data cleanup:
#    - Rows 26/27 had the "Terms Typically Offered" text run together with a
#      stray "Prerequisites:" fragment; strip that fragment out of column G.
#    - Row 42 had a "Recommended: ..." note embedded in Prerequisites (C);
#      move it out to the new Recommended column (F) and add trailing space
#      to the offered-terms text, matching the source edit.
#    - Row 62 had a "Concurrent: ..." note embedded in Prerequisites (C);
#      move it out to the new Concurrent column (E).
# ---------------------------------------------------------------------------
$ws.Cells.Item(26,7).Value = "Junior standing or History major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area D3."
$ws.Cells.Item(27,7).Value = "F, Junior standing or History major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area D3."

$ws.Cells.Item(42,3).Value = "Junior standing or History major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area D1, D2, or D3."
$ws.Cells.Item(42,6).Value = "One or more courses in GE Area B."
$ws.Cells.Item(42,7).Value = "SP "

$nbsp = [char]0x00A0
$ws.Cells.Item(59,3).Value = "HIST" + $nbsp + "303 or completion of GE Area D5, or graduate standing. Recommended GE HIST" + $nbsp + "316, HUM 310."

$ws.Cells.Item(62,3).Value = "HIST" + $nbsp + "424."
$ws.Cells.Item(62,5).Value = "EDUC" + $nbsp + "469 or EDUC" + $nbsp + "479."
$ws.Cells.Item(62,7).Value = "SP "
